$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-01 03:48:08"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
